# Weekly price update: a new weekly record is inserted at row 224
# (Fecha 2022-06-10 / serial 44722), pushing the existing rows 224-275
# down to 225-276.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 224, shifting rows 224:275 down to 225:276.
$ws.Rows("224:224").Insert()

# Populate the newly inserted row 224 with the new weekly data point.
$ws.Range("A224").Value = 4
$ws.Range("B224").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C224").Value = 'Los Lagos'
$ws.Range("D224").Value = 44722
$ws.Range("E224").Value = 10
$ws.Range("F224").Value = 100112003
$ws.Range("G224").Value = 'Ajo'
$ws.Range("H224").Value = 'Chino'
$ws.Range("I224").Value = 'Primera'
$ws.Range("J224").Value = 190
$ws.Range("K224").Value = 22000
$ws.Range("L224").Value = 22000
$ws.Range("M224").Value = 22000
$ws.Range("N224").Value = '$/caja 10 kilos'
$ws.Range("O224").Value = 'China'
$ws.Range("P224").Value = 2200
$ws.Range("Q224").Value = 10
$ws.Range("R224").Value = 'Hortaliza'
